# Auto-generated edit script: updates Leve profit calculation columns (H-N)
# across multiple worksheets to reflect refreshed market-board pricing data.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 903.6875
$ws.Range("I92").Value = 961.9
$ws.Range("J92").Value = 806.6667
$ws.Range("K92").Value = 961.9
$ws.Range("L92").Value = 806.6667
$ws.Range("M92").Value = 286.1
$ws.Range("N92").Value = -3302.6667

$ws.Range("H103").Value = 2565.1
$ws.Range("I103").Value = 2862.5
$ws.Range("J103").Value = 2366.8333
$ws.Range("K103").Value = 8587.5
$ws.Range("L103").Value = 7100.499899999999
$ws.Range("M103").Value = -8001.5
$ws.Range("N103").Value = -8272.499899999999

$ws.Range("H113").Value = 64472.75
$ws.Range("I113").Value = 144909.14
$ws.Range("K113").Value = 144909.14
$ws.Range("M113").Value = -141655.14

$ws.Range("H138").Value = 3190.5264
$ws.Range("I138").Value = 2222.2104
$ws.Range("J138").Value = 3432.6052
$ws.Range("K138").Value = 6666.6312
$ws.Range("L138").Value = 10297.8156
$ws.Range("M138").Value = -1526.6312
$ws.Range("N138").Value = -20577.8156

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 12000
$ws.Range("J23").Value = 12000
$ws.Range("L23").Value = 12000
$ws.Range("N23").Value = -12518

$ws.Range("H44").Value = 5820
$ws.Range("J44").Value = 5820
$ws.Range("L44").Value = 5820
$ws.Range("N44").Value = -6796

$ws.Range("H45").Value = 60415.41
$ws.Range("I45").Value = 92084
$ws.Range("K45").Value = 92084
$ws.Range("M45").Value = -91707

$ws.Range("H55").Value = 14675
$ws.Range("J55").Value = 15485.714
$ws.Range("L55").Value = 15485.714
$ws.Range("N55").Value = -16115.714

$ws.Range("H63").Value = 3000
$ws.Range("I63").Value = 1800
$ws.Range("J63").Value = 3400
$ws.Range("K63").Value = 1800
$ws.Range("L63").Value = 3400
$ws.Range("M63").Value = -1114
$ws.Range("N63").Value = -4772

$ws.Range("H64").Value = 42985.25
$ws.Range("J64").Value = 42985.25
$ws.Range("L64").Value = 42985.25
$ws.Range("N64").Value = -43481.25

$ws.Range("H66").Value = 3000
$ws.Range("I66").Value = 1800
$ws.Range("J66").Value = 3400
$ws.Range("K66").Value = 9000
$ws.Range("L66").Value = 17000
$ws.Range("M66").Value = -5568
$ws.Range("N66").Value = -23864

$ws.Range("H67").Value = 42985.25
$ws.Range("J67").Value = 42985.25
$ws.Range("L67").Value = 42985.25
$ws.Range("N67").Value = -44701.25

$ws.Range("H80").Value = 24585.334
$ws.Range("J80").Value = 25502.4
$ws.Range("L80").Value = 25502.4
$ws.Range("N80").Value = -27498.4

$ws.Range("H83").Value = 24585.334
$ws.Range("J83").Value = 25502.4
$ws.Range("L83").Value = 76507.20000000001
$ws.Range("N83").Value = -86491.20000000001

$ws.Range("H110").Value = 52743336
$ws.Range("I110").Value = 55673452
$ws.Range("K110").Value = 55673452
$ws.Range("M110").Value = -55671407

$ws.Range("H122").Value = 2461.04
$ws.Range("I122").Value = 2101
$ws.Range("J122").Value = 3101.111
$ws.Range("K122").Value = 6303
$ws.Range("L122").Value = 9303.332999999999
$ws.Range("M122").Value = -3853
$ws.Range("N122").Value = -14203.333

$ws.Range("H132").Value = 18719.314
$ws.Range("I132").Value = 26867.783
$ws.Range("J132").Value = 3101.4167
$ws.Range("K132").Value = 80603.349
$ws.Range("L132").Value = 9304.250100000001
$ws.Range("M132").Value = -78073.349
$ws.Range("N132").Value = -14364.2501

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 45620
$ws.Range("J62").Value = 45620
$ws.Range("L62").Value = 45620
$ws.Range("N62").Value = -46992

$ws.Range("H65").Value = 45620
$ws.Range("J65").Value = 45620
$ws.Range("L65").Value = 136860
$ws.Range("N65").Value = -143724

$ws.Range("H105").Value = 401678
$ws.Range("I105").Value = 335393
$ws.Range("J105").Value = 501105.5
$ws.Range("K105").Value = 335393
$ws.Range("L105").Value = 501105.5
$ws.Range("M105").Value = -333646
$ws.Range("N105").Value = -504599.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 800.6667
$ws.Range("I122").Value = 737.5714
$ws.Range("K122").Value = 2212.7142
$ws.Range("M122").Value = 237.2857999999997

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 1021604.8
$ws.Range("J37").Value = 1021604.8
$ws.Range("L37").Value = 3064814.4
$ws.Range("N37").Value = -3065038.4

$ws.Range("H113").Value = 812.8
$ws.Range("J113").Value = 538.26086
$ws.Range("L113").Value = 1614.78258
$ws.Range("N113").Value = -5954.78258

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3173.3635
$ws.Range("I122").Value = 2487.5
$ws.Range("K122").Value = 7462.5
$ws.Range("M122").Value = -5012.5

$ws.Range("H132").Value = 3978.85
$ws.Range("I132").Value = 2612.8572
$ws.Range("K132").Value = 7838.571599999999
$ws.Range("M132").Value = -5308.571599999999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2736
$ws.Range("I7").Value = 1860
$ws.Range("J7").Value = 4050
$ws.Range("K7").Value = 1860
$ws.Range("L7").Value = 4050
$ws.Range("M7").Value = -1748
$ws.Range("N7").Value = -4274

$ws.Range("H68").Value = 4034.2
$ws.Range("I68").Value = 2200.2856
$ws.Range("J68").Value = 8313.333000000001
$ws.Range("K68").Value = 2200.2856
$ws.Range("L68").Value = 8313.333000000001
$ws.Range("M68").Value = -1451.2856
$ws.Range("N68").Value = -9811.333000000001

$ws.Range("H71").Value = 4034.2
$ws.Range("I71").Value = 2200.2856
$ws.Range("J71").Value = 8313.333000000001
$ws.Range("K71").Value = 11001.428
$ws.Range("L71").Value = 41566.665
$ws.Range("M71").Value = -7257.428
$ws.Range("N71").Value = -49054.665

$ws.Range("H122").Value = 3126.3704
$ws.Range("I122").Value = 3000.6
$ws.Range("J122").Value = 3485.7144
$ws.Range("K122").Value = 9001.799999999999
$ws.Range("L122").Value = 10457.1432
$ws.Range("M122").Value = -6551.799999999999
$ws.Range("N122").Value = -15357.1432

$ws.Range("H126").Value = 2736
$ws.Range("I126").Value = 1860
$ws.Range("J126").Value = 4050
$ws.Range("K126").Value = 5580
$ws.Range("L126").Value = 12150
$ws.Range("M126").Value = -3110
$ws.Range("N126").Value = -17090

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1965.6
$ws.Range("I122").Value = 1965.6
$ws.Range("K122").Value = 5896.799999999999
$ws.Range("M122").Value = -3446.799999999999
